$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.201.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.85%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.797.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.29%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '

# Row 6
$ws.Range("E6").Value = '  -0.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4576'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +19.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3784'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +13.09%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.24'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.147'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.52%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07589'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.04%  '

# Row 12
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.09%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.42%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.320'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.62%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.542'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.55%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.802.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.13%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001088'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.99%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06722'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.47%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.75%  '

# Row 20
$ws.Range("E20").Value = '  -0.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.93%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.393'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.265.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.21%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.69%  '

# Row 26
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.16%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.58%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.355'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.72%  '

# Row 29
$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.008.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.43%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.50%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.235'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.64%  '

# Row 32
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.024'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.12%  '

# Row 33
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09467'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.10%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.816'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '

# Row 35
$ws.Range("B35").Value = 'Algorand'
$ws.Range("C35").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2292'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.43%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02346'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.23%  '

# Row 37
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06326'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.08%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '12.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.96%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.240'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.54%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6598'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.73%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.235'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.17%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.361'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.66%  '

# Row 43
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.481'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.00%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.90%  '

# Row 45
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.31%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.871'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.09%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6084'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.37%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.70%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.028'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.16%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07152'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.36%  '

# Row 51
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.172'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.55%  '
